# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$sheet1Updates = @{
    2  = 211
    3  = 5581
    5  = 66
    7  = 670
    8  = 655
    9  = 18
    10 = 9
    14 = 5364
    16 = 269
    17 = 231
    18 = 41
    19 = 18
    20 = 116
    21 = 4469
    22 = 226
    23 = 1181
    25 = 75
    27 = 64
    28 = 188
    30 = 152
    33 = 19
    34 = 45
    36 = 13
    40 = 46
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型": row -> new F value
$sheet4Updates = @{
    2  = 211
    4  = 5581
    6  = 66
    8  = 670
    9  = 655
    10 = 18
    11 = 9
    15 = 5364
    17 = 269
    18 = 231
    19 = 41
    20 = 18
    21 = 116
    22 = 4469
    23 = 226
    24 = 1181
    26 = 75
    28 = 64
    29 = 188
    31 = 152
    34 = 19
    35 = 45
    37 = 13
    41 = 46
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
